$wb = $excel.ActiveWorkbook

# ---- Sheet "Mapping" ----
$ws1 = $wb.Worksheets.Item("Mapping")
$ws1.Columns.Item(25).Insert()
$ws1.Cells.Item(1,25).Value = "mapping_source"
$ws1.Columns.Item(31).Insert()
$ws1.Columns.Item(31).Insert()
$ws1.Cells.Item(1,31).Value = "curation_rule"
$ws1.Cells.Item(1,32).Value = "curation_rule_text"

# ---- Sheet "MappingRegistry" ----
$ws2 = $wb.Worksheets.Item("MappingRegistry")
$ws2.Columns.Item(2).Insert()
$ws2.Columns.Item(2).Insert()
$ws2.Cells.Item(1,2).Value = "mapping_registry_title"
$ws2.Cells.Item(1,3).Value = "mapping_registry_description"

# ---- Sheet "MappingSet" ----
$ws3 = $wb.Worksheets.Item("MappingSet")
$ws3.Columns.Item(5).Insert()
$ws3.Cells.Item(1,5).Value = "mapping_set_title"
